# ScopePropertiesMiddleware tests now working.
# Update the expected JSON payloads (the middleware now flattens claims/headers
# onto the result object instead of nesting them) and fix the casing of the
# "X-User" header key used for test scenario B's Headers column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (scenario A / Expected): claims & headers are now flattened onto the
# JSON object instead of being nested under "Claims"/"Headers" arrays, and
# Host / X-HostPath fields were added.
$ws.Range("G4").Value = '{"User":"jack","name":"jack","role":"admin","Host":"localhost","hdr1":"ABC","hdr2":"DEF","X-HostPath":"localhost"}'

# Row 7 (scenario B / Expected): same flattening, plus Host field.
$ws.Range("G7").Value = '{"User":"jill","role":"user","group":"456","Host":"localhost","hdr1":"123","X-User":"jill"}'

# Row 9 (scenario C / Headers): correct the header key casing to match the
# X-User casing used elsewhere ("x-user" -> "X-User").
$ws.Range("G9").Value = "header*hdr1=123&header*X-User=jill"

# Reflect the author's final selection/cursor position in the sheet.
$ws.Range("G10").Select()
